$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg)
# for rows 2-14. Row 13 is unchanged by the diff.
$rows = @{
    2  = @{ D = 44810; J = 40; K = 12000; L = 13000; M = 12500; P = 962 }
    3  = @{ D = 44846; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    4  = @{ D = 44859; J = 30; K = 13000; L = 13000; M = 13000; P = 1000 }
    5  = @{ D = 44841; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    6  = @{ D = 44874; J = 30; K = 17000; L = 17000; M = 17000; P = 1308 }
    7  = @{ D = 44804; J = 40; K = 12000; L = 13000; M = 12500; P = 962 }
    8  = @{ D = 44868; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    9  = @{ D = 44880; J = 30; K = 17000; L = 17000; M = 17000; P = 1308 }
    10 = @{ D = 44894; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    11 = @{ D = 44797; J = 60; K = 12000; L = 13000; M = 12500; P = 962 }
    12 = @{ D = 44895; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    13 = @{ D = 44839; J = 40; K = 15000; L = 16000; M = 15500; P = 1192 }
    14 = @{ D = 44832; J = 60; K = 17000; L = 18000; M = 17500; P = 1346 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("P$r").Value = $vals.P
}
